# Fruta / hortaliza, semanal
# Insert a new data row at row 478 (pushing the existing rows 478:517 down to
# 479:518) and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 478; this shifts rows 478:517
# down to 479:518 and extends the used range to A1:T518.
$ws.Range("A478").EntireRow.Insert()

$r = 478
$ws.Cells.Item($r, 1).Value  = 10
$ws.Cells.Item($r, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item($r, 3).Value  = "La Araucanía"
$ws.Cells.Item($r, 4).Value  = 45223
$ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 5).Value  = 9
$ws.Cells.Item($r, 6).Value  = "Fruta"
$ws.Cells.Item($r, 7).Value  = 100102
$ws.Cells.Item($r, 8).Value  = "Cítricos"
$ws.Cells.Item($r, 9).Value  = 100102006
$ws.Cells.Item($r, 10).Value = "Pomelo"
$ws.Cells.Item($r, 11).Value = "Start Ruby"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 55
$ws.Cells.Item($r, 14).Value = 15000
$ws.Cells.Item($r, 15).Value = 15000
$ws.Cells.Item($r, 16).Value = 15000
$ws.Cells.Item($r, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value = 1000
$ws.Cells.Item($r, 20).Value = 15
